# Wed, Jun 10, 2020  7:05:41 PM
#
# 1) Re-style the three tables (slides 14-16) from the default
#    "Themed Style 1 - Accent 1" table style to the alternate
#    built-in table style.
# 2) Swap the deck's applied theme palette ("Integral" / Red Violet)
#    for the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

$newTableStyle = "{1E752297-52F4-4A26-BFB8-88A31B9D8BC1}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# Office Theme color scheme (RRGGBB), in ThemeColorScheme.Item() order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le 12; $k++) {
    $hex = $officeColors[$k - 1]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    $tcs.Item($k).RGB = ($b * 0x10000) + ($g * 0x100) + $r
}
